# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worksheet "Hoja1" lists, per row, a worker's document number (C),
# worker name (D), the mora period (E, format YYMM) and the mora/salary
# values (F,G). Historically only worker 73153258 (HERNANDO OLIVO
# SEMACARITT) appeared, periods 1607..1712 in descending order.
#
# This update: (1) re-sorts worker 73153258's periods into ascending
# order (1607 -> 1707 straight through rows 16-28), then (2) interleaves
# the first part ("parte 1") of a brand-new worker, 8852634 (JULIO
# HERNANDO CASTELLANO RICARDO), whose statement starts at period 1707 and
# runs to 1712, one row of the new worker following each of worker
# 73153258's remaining monthly rows (rows 29-39).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$docA = "73153258"
$nameA = "HERNANDO OLIVO SEMACARITT"
$moraA = 27578
$salA  = 781242

$docB = "8852634"
$nameB = "JULIO HERNANDO CASTELLANO RICARDO"
$moraB = 200000
$salB  = 5000000

# row => (doc, name, period, mora, salario)
$rows = @{
    16 = @($docA, $nameA, "1607", $moraA, $salA)
    17 = @($docA, $nameA, "1608", $moraA, $salA)
    18 = @($docA, $nameA, "1609", $moraA, $salA)
    19 = @($docA, $nameA, "1610", $moraA, $salA)
    20 = @($docA, $nameA, "1611", $moraA, $salA)
    21 = @($docA, $nameA, "1612", $moraA, $salA)
    22 = @($docA, $nameA, "1701", $moraA, $salA)
    23 = @($docA, $nameA, "1702", $moraA, $salA)
    24 = @($docA, $nameA, "1703", $moraA, $salA)
    25 = @($docA, $nameA, "1704", $moraA, $salA)
    26 = @($docA, $nameA, "1705", $moraA, $salA)
    27 = @($docA, $nameA, "1706", $moraA, $salA)
    28 = @($docA, $nameA, "1707", $moraA, $salA)
    29 = @($docB, $nameB, "1707", $moraB, $salB)
    30 = @($docA, $nameA, "1708", $moraA, $salA)
    31 = @($docB, $nameB, "1708", $moraB, $salB)
    32 = @($docA, $nameA, "1709", $moraA, $salA)
    33 = @($docB, $nameB, "1709", $moraB, $salB)
    34 = @($docA, $nameA, "1710", $moraA, $salA)
    35 = @($docB, $nameB, "1710", $moraB, $salB)
    36 = @($docA, $nameA, "1711", $moraA, $salA)
    37 = @($docB, $nameB, "1711", $moraB, $salB)
    38 = @($docA, $nameA, "1712", $moraA, $salA)
    39 = @($docB, $nameB, "1712", $moraB, $salB)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("C$r").Value = $vals[0]
    $ws.Range("D$r").Value = $vals[1]
    $ws.Range("E$r").Value = $vals[2]
    $ws.Range("F$r").Value = $vals[3]
    $ws.Range("G$r").Value = $vals[4]
}
